# Weekly update: insert a new week's worth of Coliflor price data
# (2 rows) at the top of the data block (row 912), shifting the
# existing historical rows down by two, and fill in the new rows'
# values. This mirrors a new "latest week" entry being prepended to
# the rolling dataset while keeping the rest of the history intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 912, pushing rows 912:951 down to 914:953.
$ws.Rows("912:913").Insert()

# Row 912 - Primera quality, new week (date serial 44939)
$ws.Cells.Item(912, 1).Value = 8
$ws.Cells.Item(912, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(912, 3).Value = "Coquimbo"
$ws.Cells.Item(912, 4).Value = 44939
$ws.Cells.Item(912, 5).Value = 4
$ws.Cells.Item(912, 6).Value = 100112008
$ws.Cells.Item(912, 7).Value = "Coliflor"
$ws.Cells.Item(912, 8).Value = "Sin especificar"
$ws.Cells.Item(912, 9).Value = "Primera"
$ws.Cells.Item(912, 10).Value = 2600
$ws.Cells.Item(912, 11).Value = 900
$ws.Cells.Item(912, 12).Value = 1000
$ws.Cells.Item(912, 13).Value = 950
$ws.Cells.Item(912, 14).Value = "$/unidad"
$ws.Cells.Item(912, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(912, 16).Value = 950
$ws.Cells.Item(912, 17).Value = 1
$ws.Cells.Item(912, 18).Value = "Hortaliza"

# Row 913 - Segunda quality, same new week
$ws.Cells.Item(913, 1).Value = 8
$ws.Cells.Item(913, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(913, 3).Value = "Coquimbo"
$ws.Cells.Item(913, 4).Value = 44939
$ws.Cells.Item(913, 5).Value = 4
$ws.Cells.Item(913, 6).Value = 100112008
$ws.Cells.Item(913, 7).Value = "Coliflor"
$ws.Cells.Item(913, 8).Value = "Sin especificar"
$ws.Cells.Item(913, 9).Value = "Segunda"
$ws.Cells.Item(913, 10).Value = 1560
$ws.Cells.Item(913, 11).Value = 700
$ws.Cells.Item(913, 12).Value = 800
$ws.Cells.Item(913, 13).Value = 750
$ws.Cells.Item(913, 14).Value = "$/unidad"
$ws.Cells.Item(913, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(913, 16).Value = 750
$ws.Cells.Item(913, 17).Value = 1
$ws.Cells.Item(913, 18).Value = "Hortaliza"
